$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at the top of the existing block (row 672),
# pushing the previous rows 672-685 down to 675-688.
$ws.Rows.Item(672).EntireRow.Insert()
$ws.Rows.Item(672).EntireRow.Insert()
$ws.Rows.Item(672).EntireRow.Insert()

# New weekly data for row 672 (Primera)
$ws.Cells.Item(672, 1).Value = 1
$ws.Cells.Item(672, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(672, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(672, 4).Value = 45239
$ws.Cells.Item(672, 5).Value = 15
$ws.Cells.Item(672, 6).Value = 100114013
$ws.Cells.Item(672, 7).Value = "Zanahoria"
$ws.Cells.Item(672, 8).Value = "Sin especificar"
$ws.Cells.Item(672, 9).Value = "Primera"
$ws.Cells.Item(672, 10).Value = 45
$ws.Cells.Item(672, 11).Value = 13000
$ws.Cells.Item(672, 12).Value = 14000
$ws.Cells.Item(672, 13).Value = 13444
$ws.Cells.Item(672, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(672, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(672, 16).Value = 538
$ws.Cells.Item(672, 17).Value = 25
$ws.Cells.Item(672, 18).Value = "Hortaliza"

# New weekly data for row 673 (Segunda)
$ws.Cells.Item(673, 1).Value = 1
$ws.Cells.Item(673, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(673, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(673, 4).Value = 45239
$ws.Cells.Item(673, 5).Value = 15
$ws.Cells.Item(673, 6).Value = 100114013
$ws.Cells.Item(673, 7).Value = "Zanahoria"
$ws.Cells.Item(673, 8).Value = "Sin especificar"
$ws.Cells.Item(673, 9).Value = "Segunda"
$ws.Cells.Item(673, 10).Value = 15
$ws.Cells.Item(673, 11).Value = 11000
$ws.Cells.Item(673, 12).Value = 12000
$ws.Cells.Item(673, 13).Value = 11667
$ws.Cells.Item(673, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(673, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(673, 16).Value = 467
$ws.Cells.Item(673, 17).Value = 25
$ws.Cells.Item(673, 18).Value = "Hortaliza"

# New weekly data for row 674 (Tercera)
$ws.Cells.Item(674, 1).Value = 1
$ws.Cells.Item(674, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(674, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(674, 4).Value = 45239
$ws.Cells.Item(674, 5).Value = 15
$ws.Cells.Item(674, 6).Value = 100114013
$ws.Cells.Item(674, 7).Value = "Zanahoria"
$ws.Cells.Item(674, 8).Value = "Sin especificar"
$ws.Cells.Item(674, 9).Value = "Tercera"
$ws.Cells.Item(674, 10).Value = 20
$ws.Cells.Item(674, 11).Value = 9000
$ws.Cells.Item(674, 12).Value = 10000
$ws.Cells.Item(674, 13).Value = 9750
$ws.Cells.Item(674, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(674, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(674, 16).Value = 390
$ws.Cells.Item(674, 17).Value = 25
$ws.Cells.Item(674, 18).Value = "Hortaliza"
